$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.129.72'
$ws.Range("E2").Value = '  -3.41%  '

$ws.Range("D3").Value = '1.849.45'
$ws.Range("E3").Value = '  -2.27%  '

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = "'0.7052"
$ws.Range("E5").Value = '  -5.52%  '

$ws.Range("D6").Value = "'238.23"
$ws.Range("E6").Value = '  -1.88%  '

$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = "'0.3049"
$ws.Range("E8").Value = '  -3.90%  '

$ws.Range("D9").Value = "'0.07511"
$ws.Range("E9").Value = '  +3.63%  '

$ws.Range("D10").Value = "'23.34"
$ws.Range("E10").Value = '  -6.73%  '

$ws.Range("D11").Value = "'0.08127"
$ws.Range("E11").Value = '  -2.84%  '

$ws.Range("D12").Value = '1.864.46'
$ws.Range("E12").Value = '  -1.69%  '

$ws.Range("D13").Value = "'0.7248"
$ws.Range("E13").Value = '  -5.01%  '

$ws.Range("D14").Value = "'5.218"
$ws.Range("E14").Value = '  -4.34%  '

$ws.Range("D15").Value = "'89.12"
$ws.Range("E15").Value = '  -4.26%  '

$ws.Range("D16").Value = '29.370.69'
$ws.Range("E16").Value = '  -2.76%  '

$ws.Range("D17").Value = "'5.790"
$ws.Range("E17").Value = '  -6.48%  '

$ws.Range("D18").Value = "'239.01"
$ws.Range("E18").Value = '  -4.90%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = "'0.000007673"
$ws.Range("E19").Value = '  -2.56%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = "'13.07"
$ws.Range("E20").Value = '  -4.48%  '

$ws.Range("E21").Value = '  +0.12%  '

$ws.Range("D22").Value = '2.139.92'
$ws.Range("E22").Value = '  -0.78%  '

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").Value = "'7.571"
$ws.Range("E24").Value = '  -5.78%  '

$ws.Range("D25").Value = "'0.1465"
$ws.Range("E25").Value = '  -7.95%  '

$ws.Range("D26").Value = "'8.990"
$ws.Range("E26").Value = '  -3.56%  '

$ws.Range("D27").Value = "'161.62"
$ws.Range("E27").Value = '  -1.59%  '

$ws.Range("D28").Value = "'18.01"
$ws.Range("E28").Value = '  -4.30%  '

$ws.Range("D29").Value = "'1.937"
$ws.Range("E29").Value = '  -6.69%  '

$ws.Range("D30").Value = "'1.384"
$ws.Range("E30").Value = '  -6.27%  '

$ws.Range("D31").Value = "'4.568"

$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("D33").Value = "'4.002"
$ws.Range("E33").Value = '  -5.44%  '

$ws.Range("D34").Value = "'0.05162"
$ws.Range("E34").Value = '  -4.63%  '

$ws.Range("D35").Value = "'1.186"
$ws.Range("E35").Value = '  -5.73%  '

$ws.Range("D36").Value = "'1.041"
$ws.Range("E36").Value = '  +4.76%  '

$ws.Range("D37").Value = "'0.7030"
$ws.Range("E37").Value = '  -8.50%  '

$ws.Range("D38").Value = "'2.643"
$ws.Range("E38").Value = '  -2.65%  '

$ws.Range("E39").Value = '  -5.57%  '

$ws.Range("D40").Value = "'2.676"
$ws.Range("E40").Value = '  -3.45%  '

$ws.Range("D41").Value = "'0.9487"
$ws.Range("E41").Value = '  +8.98%  '

$ws.Range("D42").Value = "'6.010"
$ws.Range("E42").Value = '  -1.55%  '

$ws.Range("D43").Value = '1.076.68'
$ws.Range("E43").Value = '  -2.59%  '

$ws.Range("D44").Value = "'0.4303"
$ws.Range("E44").Value = '  -6.18%  '

$ws.Range("D45").Value = "'70.05"
$ws.Range("E45").Value = '  -4.22%  '

$ws.Range("D46").Value = "'0.9996"
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").Value = "'102.24"
$ws.Range("E47").Value = '  -2.17%  '

$ws.Range("D48").Value = "'1.752"
$ws.Range("E48").Value = '  -6.53%  '

$ws.Range("D49").Value = '2.011.28'
$ws.Range("E49").Value = '  -2.22%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'9.241"
$ws.Range("E50").Value = '  -4.20%  '

$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = "'7.062"
$ws.Range("E51").Value = '  -7.54%  '
